$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "key"
$ws.Range("B1").Value = "iri"
$ws.Range("C1").Value = "annotation"

$ws.Rows.Item(1).RowHeight = 19.5
$ws.Rows.Item(2).RowHeight = 19.5
